$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.239.68'
$ws.Range("E2").Value = '  +1.38%  '

$ws.Range("D3").Value = '1.653.68'
$ws.Range("E3").Value = '  +0.35%  '

$ws.Range("E4").Value = '  -0.84%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.80'
$ws.Range("E5").Value = '  +1.28%  '

$ws.Range("E6").Value = '  -0.30%  '

$ws.Range("E7").Value = '  -0.93%  '

$ws.Range("E8").Value = '  +0.65%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0627'
$ws.Range("E9").Value = '  -0.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.61'
$ws.Range("E10").Value = '  +1.87%  '

$ws.Range("E11").Value = '  +0.41%  '

$ws.Range("D12").Value = '1.882.88'
$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("D13").Value = '1.655.36'
$ws.Range("E13").Value = '  +0.31%  '

$ws.Range("E14").Value = '  +0.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.532'
$ws.Range("E15").Value = '  +0.41%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.06'
$ws.Range("E16").Value = '  +1.86%  '

$ws.Range("D17").Value = '27.184.41'
$ws.Range("E17").Value = '  +1.14%  '

$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("E18").Value = '  +0.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '221.41'
$ws.Range("E19").Value = '  +3.06%  '

$ws.Range("E20").Value = '  -0.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.76'
$ws.Range("E21").Value = '  +7.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.43'
$ws.Range("E22").Value = '  +0.78%  '

$ws.Range("E23").Value = '  -2.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.28'
$ws.Range("E24").Value = '  -0.88%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.34'
$ws.Range("E25").Value = '  -0.11%  '

$ws.Range("E26").Value = '  -0.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.36'
$ws.Range("E27").Value = '  +2.06%  '

$ws.Range("E28").Value = '  +0.35%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.97'
$ws.Range("E29").Value = '  +1.81%  '

$ws.Range("E30").Value = '  +1.56%  '

$ws.Range("E31").Value = '  +0.78%  '

$ws.Range("E32").Value = '  +0.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.01'
$ws.Range("E33").Value = '  -0.26%  '

$ws.Range("E34").Value = '  +2.75%  '

$ws.Range("D35").Value = '1.267.90'
$ws.Range("E35").Value = '  -2.29%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.44'
$ws.Range("E36").Value = '  -0.17%  '

$ws.Range("E37").Value = '  -1.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.538'
$ws.Range("E38").Value = '  +0.44%  '

$ws.Range("E39").Value = '  +0.46%  '

$ws.Range("E40").Value = '  -0.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.810'
$ws.Range("E41").Value = '  +0.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.38'
$ws.Range("E42").Value = '  +0.82%  '

$ws.Range("D43").Value = '1.793.01'
$ws.Range("E43").Value = '  +0.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.04'
$ws.Range("E44").Value = '  +0.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.67'
$ws.Range("E45").Value = '  +0.64%  '

$ws.Range("E46").Value = '  -7.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.61'
$ws.Range("E47").Value = '  +0.15%  '

$ws.Range("E48").Value = '  -0.68%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0976'
$ws.Range("E49").Value = '  +0.51%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.64'
$ws.Range("E50").Value = '  -0.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.405'
$ws.Range("E51").Value = '  -0.57%  '
